$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching style of other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill H2:H9 with 0
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
